# Fruta / hortaliza, semanal
# Weekly data refresh: insert a new daily-price record as row 227 (pushing
# every existing record down by one row) for Apio / Vega Monumental
# Concepción, dated 2022-08-11 (serial 44784).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 227..318 down to 228..319 by inserting a fresh row at 227.
$ws.Rows.Item(227).EntireRow.Insert()

# Populate the newly inserted row 227 with the new weekly observation.
# Most fields mirror the record that used to sit at row 227 (same market,
# region, quality grade, volume, unit, etc.) — only the date and the
# min/max/avg price + $/Kg columns are new.
$ws.Range("A227").Value = 11
$ws.Range("B227").Value = "Vega Monumental Concepción"
$ws.Range("C227").Value = "Bíobío"
$ws.Range("D227").Value = 44784
$ws.Range("E227").Value = 8
$ws.Range("F227").Value = 100112017
$ws.Range("G227").Value = "Apio"
$ws.Range("H227").Value = "Americana (o)"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 100
$ws.Range("K227").Value = 8000
$ws.Range("L227").Value = 8500
$ws.Range("M227").Value = 8250
$ws.Range("N227").Value = "$/docena de matas"
$ws.Range("O227").Value = "Región de Coquimbo"
$ws.Range("P227").Value = 1375
$ws.Range("Q227").Value = 6
$ws.Range("R227").Value = "Hortaliza"
